# Update database and shift reporting periods (read_price algorithm change):
# Each period column shifts one column to the left (D<-E<-F<-G<-H) and the
# newest period / publish-date / figures are written into column H.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8: period headers -------------------------------------------------
$ws.Range("D8").Value = "12 ماهه منتهی به 1397/12"
$ws.Range("E8").Value = "12 ماهه منتهی به 1398/12"
$ws.Range("F8").Value = "12 ماهه منتهی به 1399/12"
$ws.Range("G8").Value = "12 ماهه منتهی به 1400/12"
$ws.Range("H8").Value = "12 ماهه منتهی به 1401/12"

# --- Row 9: publish dates ----------------------------------------------------
$ws.Range("D9").Value = "1399-04-16 (8)"
$ws.Range("E9").Value = "1400-04-20 (8)"
$ws.Range("F9").Value = "1401-04-11 (8)"
$ws.Range("G9").Value = "1402-02-28 (7)"
# "1402-02-28" alone looks exactly like a calendar date to Excel's literal
# parser and would otherwise be auto-converted to a date serial number, so
# force the cell to Text format first, then restore the original cell
# formatting (border/fill/alignment) by copying it over from its neighbour.
$ws.Range("H9").NumberFormat = "@"
$ws.Range("H9").Value = "1402-02-28"
$ws.Range("G9").Copy()
$ws.Range("H9").PasteSpecial(-4122)

# --- Row 11: فروش ------------------------------------------------------------
$ws.Range("D11").Value = 32264
$ws.Range("E11").Value = 44543
$ws.Range("F11").Value = 49786
$ws.Range("G11").Value = 54796
$ws.Range("H11").Value = 58031

# --- Row 12: بهای تمام شده کالای فروش رفته -----------------------------------
$ws.Range("D12").Value = -25460
$ws.Range("E12").Value = -35067
$ws.Range("F12").Value = -31249
$ws.Range("G12").Value = -47481
$ws.Range("H12").Value = -50720

# --- Row 13: سود (زیان) ناخالص ------------------------------------------------
$ws.Range("D13").Value = 6803
$ws.Range("E13").Value = 9477
$ws.Range("F13").Value = 18537
$ws.Range("G13").Value = 7316
$ws.Range("H13").Value = 7311

# --- Row 14: هزینه های عمومی, اداری و تشکیلاتی --------------------------------
$ws.Range("D14").Value = -1349
$ws.Range("E14").Value = -1440
$ws.Range("F14").Value = -1148
$ws.Range("G14").Value = -1403
$ws.Range("H14").Value = -1721

# --- Row 16: خالص سایر درامدها (هزینه ها) ی عملیاتی ---------------------------
$ws.Range("D16").Value = -559
$ws.Range("E16").Value = -17
$ws.Range("F16").Value = 160
$ws.Range("G16").Value = 308
$ws.Range("H16").Value = 442

# --- Row 17: سود (زیان) عملیاتی ------------------------------------------------
$ws.Range("D17").Value = 4895
$ws.Range("E17").Value = 8019
$ws.Range("F17").Value = 17549
$ws.Range("G17").Value = 6221
$ws.Range("H17").Value = 6032

# --- Row 18: هزینه های مالی -----------------------------------------------------
$ws.Range("D18").Value = -603
$ws.Range("E18").Value = -306
$ws.Range("F18").Value = -315
$ws.Range("G18").Value = -560
$ws.Range("H18").Value = -1165

# --- Row 19: خالص سایر درامدها و هزینه های غیرعملیاتی ---------------------------
$ws.Range("D19").Value = 102
$ws.Range("E19").Value = 718
$ws.Range("F19").Value = 131
$ws.Range("G19").Value = 170
$ws.Range("H19").Value = -83

# --- Row 20: سود (زیان) خالص عملیات در حال تداوم قبل از مالیات ------------------
$ws.Range("D20").Value = 4394
$ws.Range("E20").Value = 8431
$ws.Range("F20").Value = 17364
$ws.Range("G20").Value = 5832
$ws.Range("H20").Value = 4784

# --- Row 21: مالیات (was "-" in column D, now a real number because of the shift) --
$ws.Range("D21").Value = -1057
$ws.Range("E21").Value = -1702
$ws.Range("F21").Value = -1718
$ws.Range("G21").Value = -1183
$ws.Range("H21").Value = -654

# --- Row 22: سود (زیان) خالص عملیات در حال تداوم --------------------------------
$ws.Range("D22").Value = 3337
$ws.Range("E22").Value = 6729
$ws.Range("F22").Value = 15646
$ws.Range("G22").Value = 4649
$ws.Range("H22").Value = 4129

# --- Row 24: سود (زیان) خالص -----------------------------------------------------
$ws.Range("D24").Value = 3337
$ws.Range("E24").Value = 6729
$ws.Range("F24").Value = 15646
$ws.Range("G24").Value = 4649
$ws.Range("H24").Value = 4129

# --- Row 26: سرمایه ----------------------------------------------------------------
$ws.Range("D26").Value = 2832
$ws.Range("E26").Value = 2233
$ws.Range("F26").Value = 57496
$ws.Range("G26").Value = 49270
$ws.Range("H26").Value = 36839
